$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old SUM formula lived in D10; clear it first (it will move down to D13
# once the new item rows are inserted below the existing data).
$ws.Range("D10").ClearContents()

# Update existing D4 value (2nd "3D print" row price)
$ws.Range("D4").Value = 69.69

# Add new purchased item rows 10-12, matching the table's centered style
$newRange = $ws.Range("B10:D12")
$newRange.HorizontalAlignment = -4108
$newRange.VerticalAlignment = -4108

$ws.Range("B10").Value = 8
$ws.Range("C10").Value = "XT60"
$ws.Range("D10").Value = 3.6

$ws.Range("B11").Value = 9
$ws.Range("C11").Value = "EC5"
$ws.Range("D11").Value = 3.48

$ws.Range("B12").Value = 10
$ws.Range("C12").Value = "Pro Mini"
$ws.Range("D12").Value = 6.45

# Move the SUM formula down to row 13, covering the new rows
$ws.Range("D13").Formula = "=SUM(D3:D12)"
$ws.Range("D13").HorizontalAlignment = -4108
$ws.Range("D13").VerticalAlignment = -4108

# Update selection to match the target state
$ws.Range("G15").Select()
